# "cleaned and optimized for SKU number naming"
#
# The SKU column (A) is converted from text SKU codes to plain sequential
# numbers, and the "Associated Patch" column (B) is renamed "PATCH" and
# trimmed down to just a couple of short patch-letter entries ("A"/"B"),
# with two new rows added at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B cleanup -------------------------------------------------
# Rows 2, 4 and 5 currently hold leftover "Associated Patch" text values
# that are no longer needed. Clear them out completely (they still have
# the sheet's default style at this point, so ClearContents drops the
# cells entirely rather than leaving an empty styled cell behind).
$ws.Range("B2").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()

# --- Column A: SKU codes -> plain sequential numbers -------------------
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 9

# --- Column B: a couple of short patch codes, new header last ----------
# (set before the header so the shared-string table fills in the same
# order as the authored workbook: SKU, A, B, PATCH)
$ws.Range("B3").Value = "A"
$ws.Range("B8").Value = "B"
$ws.Range("B1").Value = "PATCH"

# --- Formatting: center-align the populated cells -----------------------
$ws.Range("A1:A9").HorizontalAlignment = -4108
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B8").HorizontalAlignment = -4108

# --- Selection, matching the saved workbook state -----------------------
$ws.Range("I12").Select()
